$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.646.22'
$ws.Range("E2").Value = '  +1.12%  '
$ws.Range("D3").Value = '1.632.57'
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("E6").Value = '  +2.73%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("E8").Value = '  +1.69%  '
$ws.Range("E9").Value = '  +0.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.15'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0844'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.58%  '
$ws.Range("D12").Value = '1.860.37'
$ws.Range("E12").Value = '  +0.82%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.08'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.49%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.590.94'
$ws.Range("E14").Value = '  -2.72%  '
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("D16").Value = '26.635.42'
$ws.Range("E16").Value = '  +1.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.21'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.00%  '
$ws.Range("E18").Value = '  +1.72%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '218.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.92%  '
$ws.Range("E20").Value = '  +0.15%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.17'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.16%  '
$ws.Range("E23").Value = '  +0.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.94'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.35%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("E27").Value = '  +0.65%  '
$ws.Range("E28").Value = '  +4.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.49'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.04%  '
$ws.Range("E30").Value = '  -1.76%  '
$ws.Range("E31").Value = '  -0.23%  '
$ws.Range("E32").Value = '  +3.42%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.97'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.38%  '
$ws.Range("E34").Value = '  +0.35%  '
$ws.Range("E35").Value = '  +0.55%  '
$ws.Range("D36").Value = '1.208.68'
$ws.Range("E36").Value = '  +2.70%  '
$ws.Range("E37").Value = '  +5.20%  '
$ws.Range("E38").Value = '  -0.28%  '
$ws.Range("E39").Value = '  +0.22%  '
$ws.Range("E40").Value = '  +0.56%  '
$ws.Range("E41").Value = '  -1.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.41'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.54%  '
$ws.Range("E43").Value = '  +0.30%  '
$ws.Range("D44").Value = '1.772.06'
$ws.Range("E44").Value = '  +0.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.84'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.54'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.69'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.58%  '
$ws.Range("E48").Value = '  +1.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.58'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.01%  '
$ws.Range("E50").Value = '  -0.12%  '
$ws.Range("E51").Value = '  +0.25%  '
